$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds values stored as text (inlineStr) in the
# original workbook. Setting NumberFormat to text ("@") before assigning
# the Value keeps Excel from re-interpreting the strings as numbers, which
# would otherwise drop significant trailing/leading zeros.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "270.02"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.87"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.362"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.644"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.697"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.371"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8340"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01378"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08405"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03485"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03148"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09324"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.889"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001717"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04839"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006220"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.003535"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001499"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.741"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.344"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04677"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006993"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1171"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003451"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01245"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006256"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.8793"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.08588"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01239"
